# Update the cryptocurrency price (column D) and 1h volume change (column E)
# values on the active worksheet, matching the latest scrape from the
# GitHub Actions "Updated cryptos list" job.
#
# Column D values are written as plain text (matching the source data,
# which stores them as inline strings rather than numbers), so each D
# cell is briefly formatted as Text before the assignment and then has
# its formatting cleared again to avoid leaving a stray number format
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.590.59'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +4.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.273.16'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.88%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.89'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +6.54%  '
$ws.Range("E7").Value = '  +0.86%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.66'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.54%  '
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.89'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.885'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.622.24'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.61'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.275.30'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.447.38'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.96'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.34%  '
$ws.Range("E20").Value = '  +4.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.56'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.56'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.21'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '239.98'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("E25").Value = '  +3.94%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.18'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.58'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +12.80%  '
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.49'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.76%  '
$ws.Range("E31").Value = '  +0.87%  '
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '161.74'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +4.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.78'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.117'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +10.31%  '
$ws.Range("E36").Value = '  +4.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.15'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.94'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.47'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.64'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +25.27%  '
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.771.67'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.83%  '
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '87.34'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.47'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '60.72'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.80'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.03%  '
$ws.Range("E50").Value = '  +7.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '104.56'
$ws.Range("D51").ClearFormats()
